$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 47; existing rows 47-62 shift down to 48-63.
$ws.Rows.Item(47).Insert()

# Populate the new row 47 with the new weekly record.
$ws.Range("A47").Value = 10
$ws.Range("B47").Value = "Vega Modelo de Temuco"
$ws.Range("C47").Value = "La Araucanía"
$ws.Range("D47").Value = 44523
$ws.Range("E47").Value = 9
$ws.Range("F47").Value = 300000000
$ws.Range("G47").Value = "Espárragos"
$ws.Range("H47").Value = "Sin especificar"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 130
$ws.Range("K47").Value = 1300
$ws.Range("L47").Value = 1400
$ws.Range("M47").Value = 1338
$ws.Range("N47").Value = "`$/kilo"
$ws.Range("O47").Value = "Región del Maule"
$ws.Range("P47").Value = 1338
$ws.Range("Q47").Value = 1
$ws.Range("R47").Value = "Hortaliza"

# Match the number format/style used by the other rows in column D.
$ws.Range("D47").NumberFormat = $ws.Range("D46").NumberFormat
